# recognizer_config.xlsx update:
#  - move "auto_receipt" row to the top of the data block (row 2)
#  - append 7 new OCR recognizer types (taxi/train/household/passport/driver/vehicle license/vehicle cert)
#  - refresh hyperlinks for every RecognizerUrl cell
#  - widen columns A/B to fit the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop every existing hyperlink up front - we rebuild them all once the
# table rows are in their final places (Range.Hyperlinks.Delete() clears
# the whole sheet's collection in this host).
$ws.Hyperlinks.Delete()

# --- header row -------------------------------------------------------
$ws.Range("A1").Value = "RecognizerName"
$ws.Range("B1").Value = "RecognizerUrl"
$ws.Range("C1").Value = "AppKey"
$ws.Range("D1").Value = "AppSecret"

# --- data rows that re-use already-existing strings (just re-ordered) --
$ws.Range("A2").Value = "auto_receipt"
$ws.Range("B2").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/general_receipt_recog"

$ws.Range("A3").Value = "vat_invoice"
$ws.Range("B3").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/vat_invoice"

$ws.Range("A4").Value = "bank_card"
$ws.Range("B4").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/bank_card"

$ws.Range("A5").Value = "business_license"
$ws.Range("B5").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/business_license"

$ws.Range("A6").Value = "id_card"
$ws.Range("B6").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/id_card"

# --- brand-new rows: names (column A) entered first, then urls (column B),
# matching the order the new shared strings show up in the saved file -----
$ws.Range("A7").Value = "taxi_receipt"
$ws.Range("A8").Value = "train_ticket"
$ws.Range("A9").Value = "household_register"
$ws.Range("A10").Value = "passport"
$ws.Range("A11").Value = "driver_license"
$ws.Range("A12").Value = "vehicle_license"
$ws.Range("A13").Value = "vehicle_certificate"

$ws.Range("B7").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/taxi_invoice"
$ws.Range("B8").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/train_ticket"
$ws.Range("B9").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/household_register"
$ws.Range("B10").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/passport"
$ws.Range("B11").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/driver_license"
$ws.Range("B12").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/vehicle_license"
$ws.Range("B13").Value = "https://ocr-api.ccint.com/cci_ai/service/v1/vehicle_inspection_certificate"

# the old D2 "stray" empty-but-styled cell now lives under the row it
# travelled with (vat_invoice, now row 3)
$ws.Range("D2").Clear()
$ws.Range("D3").HorizontalAlignment = 1

# --- hyperlinks (rebuilt in final row order) ---------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://ocr-api.ccint.com/cci_ai/service/v1/general_receipt_recog")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://ocr-api.ccint.com/cci_ai/service/v1/vat_invoice")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://ocr-api.ccint.com/cci_ai/service/v1/bank_card")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://ocr-api.ccint.com/cci_ai/service/v1/business_license")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://ocr-api.ccint.com/cci_ai/service/v1/id_card")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://ocr-api.ccint.com/cci_ai/service/v1/taxi_invoice")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://ocr-api.ccint.com/cci_ai/service/v1/train_ticket")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://ocr-api.ccint.com/cci_ai/service/v1/household_register")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://ocr-api.ccint.com/cci_ai/service/v1/passport")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://ocr-api.ccint.com/cci_ai/service/v1/driver_license")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://ocr-api.ccint.com/cci_ai/service/v1/vehicle_license")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://ocr-api.ccint.com/cci_ai/service/v1/vehicle_inspection_certificate")

# re-apply the Hyperlink cell style (blue/underline) to every link cell -
# Hyperlinks.Add alone leaves the default font in this host
$ws.Range("B2:B13").Style = "Hyperlink"

# --- column widths (auto-fit to new, longer content) -------------------
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 58.16666666666667

# --- selection, matching the saved workbook state -----------------------
$ws.Range("A17").Select() | Out-Null
